# Auto-generated script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.757.10"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.566.20"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "206.20"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.96%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.488"
$c.ClearFormats()
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.15%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0863"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.787.50"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.561.16"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -2.21%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.23%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "61.53"
$c.ClearFormats()
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "26.782.90"
$ws.Range("E17").Value = "  -2.35%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "214.59"
$c.ClearFormats()
$ws.Range("E18").Value = "  +1.03%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.36"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "0.0₃0677"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  -1.57%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.74"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.12%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.74"
$c.ClearFormats()
$ws.Range("E26").Value = "  +1.22%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "14.90"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0463"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.11"
$c.ClearFormats()
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.386.84"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("E39").Value = "  -0.44%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.818"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.09%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.992"
$c.ClearFormats()
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +1.08%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "63.23"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "1.701.24"
$ws.Range("E47").Value = "  +0.39%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "85.50"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "0.0₇0993"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  -0.73%  "
